# Updates cryptocurrency price/volume figures on Sheet1 (rows 2-51)
# to match the latest scrape, per the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.711.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "'2.083.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'228.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'60.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").Value = "'0.0837"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "'2.394.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "'14.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").Value = "'21.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").Value = "'0.794"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").Value = "'5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'2.090.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "'38.650.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("D19").Value = "'71.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "'226.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'170.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "'9.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  +6.94%  "
$ws.Range("E29").Value = "  +11.61%  "
$ws.Range("D30").Value = "'19.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "'2.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.68%  "
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("E34").Value = "  +4.61%  "
$ws.Range("D35").Value = "'0.0613"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "'18.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").Value = "'1.538.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'100.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").Value = "'7.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.09%  "
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'2.286.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
